# Update cryptos list with latest scraped price/volume data
# (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''26.732.70'
$ws.Range("E2").Value = '  +0.05%  '
# Row 3
$ws.Range("D3").Value = '''1.640.17'
$ws.Range("E3").Value = '  -0.25%  '
# Row 4
$ws.Range("E4").Value = '  +0.23%  '
# Row 5
$ws.Range("D5").Value = '''218.31'
$ws.Range("E5").Value = '  +1.10%  '
# Row 6
$ws.Range("E6").Value = '  -0.60%  '
# Row 7
$ws.Range("E7").Value = '  +0.12%  '
# Row 8
$ws.Range("E8").Value = '  -0.28%  '
# Row 9
$ws.Range("E9").Value = '  -0.49%  '
# Row 10
$ws.Range("D10").Value = '''19.06'
$ws.Range("E10").Value = '  -0.55%  '
# Row 11
$ws.Range("D11").Value = '''0.0845'
$ws.Range("E11").Value = '  +0.49%  '
# Row 12
$ws.Range("D12").Value = '''1.867.07'
$ws.Range("E12").Value = '  -0.39%  '
# Row 13
$ws.Range("D13").Value = '''1.646.28'
$ws.Range("E13").Value = '  -0.05%  '
# Row 14
$ws.Range("E14").Value = '  -1.04%  '
# Row 15
$ws.Range("E15").Value = '  -1.22%  '
# Row 16
$ws.Range("D16").Value = '''64.58'
$ws.Range("E16").Value = '  -0.91%  '
# Row 17
$ws.Range("D17").Value = '''26.702.69'
$ws.Range("E17").Value = '  -0.15%  '
# Row 18
$ws.Range("E18").Value = '  -2.13%  '
# Row 19
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '''1.01'
$ws.Range("E19").Value = '  +0.25%  '
# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''211.37'
$ws.Range("E20").Value = '  -2.98%  '
# Row 21
$ws.Range("E21").Value = '  -0.34%  '
# Row 22
$ws.Range("E22").Value = '  -0.96%  '
# Row 23
$ws.Range("D23").Value = '''2.32'
$ws.Range("E23").Value = '  -5.02%  '
# Row 24
$ws.Range("E24").Value = '  -2.19%  '
# Row 25
$ws.Range("D25").Value = '''147.00'
$ws.Range("E25").Value = '  +0.56%  '
# Row 26
$ws.Range("E26").Value = '  +0.17%  '
# Row 27
$ws.Range("E27").Value = '  -1.70%  '
# Row 28
$ws.Range("E28").Value = '  -0.92%  '
# Row 29
$ws.Range("D29").Value = '''15.57'
$ws.Range("E29").Value = '  -0.85%  '
# Row 30
$ws.Range("D30").Value = '''0.0501'
$ws.Range("E30").Value = '  -3.41%  '
# Row 31
$ws.Range("D31").Value = '''1.19'
$ws.Range("E31").Value = '  +0.94%  '
# Row 32
$ws.Range("E32").Value = '  +0.35%  '
# Row 33
$ws.Range("D33").Value = '''2.98'
$ws.Range("E33").Value = '  -0.41%  '
# Row 34
$ws.Range("D34").Value = '''1.272.70'
$ws.Range("E34").Value = '  -0.34%  '
# Row 35
$ws.Range("E35").Value = '  -0.73%  '
# Row 36
$ws.Range("D36").Value = '''2.45'
$ws.Range("E36").Value = '  +0.15%  '
# Row 37
$ws.Range("D37").Value = '''0.0175'
$ws.Range("E37").Value = '  -1.80%  '
# Row 38
$ws.Range("E38").Value = '  -1.73%  '
# Row 39
$ws.Range("D39").Value = '''0.805'
$ws.Range("E39").Value = '  -2.67%  '
# Row 40
$ws.Range("E40").Value = '  +0.20%  '
# Row 41
$ws.Range("D41").Value = '''0.803'
$ws.Range("E41").Value = '  -1.11%  '
# Row 42
$ws.Range("E42").Value = '  -2.19%  '
# Row 43
$ws.Range("D43").Value = '''1.778.06'
$ws.Range("E43").Value = '  -0.52%  '
# Row 44
$ws.Range("E44").Value = '  -3.39%  '
# Row 45
$ws.Range("D45").Value = '''91.39'
$ws.Range("E45").Value = '  -0.47%  '
# Row 46
$ws.Range("D46").Value = '''60.16'
# Row 47
$ws.Range("E47").Value = '  -1.28%  '
# Row 48
$ws.Range("D48").Value = '''0.0519'
$ws.Range("E48").Value = '  +0.80%  '
# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.0961'
$ws.Range("E49").Value = '  -0.80%  '
# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''7.51'
$ws.Range("E50").Value = '  -2.76%  '
# Row 51
$ws.Range("E51").Value = '  -0.04%  '
